$d = $word.ActiveDocument

# 1) Replace the thesis title text.
$d.Content.Find.Execute(
    "Разработка информационного портала организации", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Современные технологии интернет-программирования", 2) | Out-Null

# 2) Replace "Кофе центр" with "Разработка информационного портала органи" +
#    "зации", split into two runs with the (relocated) "_GoBack" bookmark
#    sitting between them — mirroring where Word happened to leave the
#    caret when the document was last saved.
$target = $d.Content
$target.Find.Execute("Кофе центр") | Out-Null
$start = $target.Start
$end = $target.End
$mid = $start + 4

# Temporary bookmarks pin the outer edges of "Кофе центр" so later text
# edits don't get coalesced into the neighbouring "«" / "»" runs.
$d.Bookmarks.Add("ZZZ_LEFT", $d.Range($start, $start)) | Out-Null
$d.Bookmarks.Add("ZZZ_RIGHT", $d.Range($end, $end)) | Out-Null

# Re-home "_GoBack" at the future split point (between "органи" and
# "зации"); this also removes it from wherever it previously sat.
$d.Bookmarks.Add("_GoBack", $d.Range($mid, $mid)) | Out-Null

$leftStart = $d.Bookmarks("ZZZ_LEFT").Range.Start
$splitPoint = $d.Bookmarks("_GoBack").Range.Start

$left = $d.Range($leftStart, $splitPoint)
$left.Text = "Разработка информационного портала органи"

$splitPoint = $d.Bookmarks("_GoBack").Range.Start
$rightEnd = $d.Bookmarks("ZZZ_RIGHT").Range.Start

$right = $d.Range($splitPoint, $rightEnd)
$right.Text = "зации"

$d.Bookmarks("ZZZ_LEFT").Delete()
$d.Bookmarks("ZZZ_RIGHT").Delete()
